# Auto-generated edit script: updates columns H-N for specific rows
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# per the scheduled-runner price/profit refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Cells.Item(21, 8).Value = 0  # H21: 9500 -> 0
$ws.Cells.Item(21, 9).Value = 0  # I21: 6000 -> 0
$ws.Cells.Item(21, 10).Value = 0  # J21: 10000 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 6000 -> 0
$ws.Cells.Item(21, 12).Value = 0  # L21: 10000 -> 0
$ws.Cells.Item(21, 13).ClearContents()  # M21: -5532 -> (blank)
$ws.Cells.Item(21, 14).ClearContents()  # N21: -10936 -> (blank)

# Row 23
$ws.Cells.Item(23, 8).Value = 0  # H23: 9500 -> 0
$ws.Cells.Item(23, 9).Value = 0  # I23: 6000 -> 0
$ws.Cells.Item(23, 10).Value = 0  # J23: 10000 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 6000 -> 0
$ws.Cells.Item(23, 12).Value = 0  # L23: 10000 -> 0
$ws.Cells.Item(23, 13).ClearContents()  # M23: -5766 -> (blank)
$ws.Cells.Item(23, 14).ClearContents()  # N23: -10468 -> (blank)

# Row 29
$ws.Cells.Item(29, 8).Value = 90  # H29: 838 -> 90
$ws.Cells.Item(29, 9).Value = 90  # I29: 95 -> 90
$ws.Cells.Item(29, 10).Value = 0  # J29: 1333.3334 -> 0
$ws.Cells.Item(29, 11).Value = 270  # K29: 285 -> 270
$ws.Cells.Item(29, 12).Value = 0  # L29: 4000.0002 -> 0
$ws.Cells.Item(29, 13).Value = 11  # M29: -4 -> 11
$ws.Cells.Item(29, 14).ClearContents()  # N29: -4562.0002 -> (blank)

# Row 40
$ws.Cells.Item(40, 8).Value = 6043907  # H40: 4835385.5 -> 6043907
$ws.Cells.Item(40, 9).Value = 8930286  # I40: 6251590 -> 8930286
$ws.Cells.Item(40, 11).Value = 8930286  # K40: 6251590 -> 8930286
$ws.Cells.Item(40, 13).Value = -8930111  # M40: -6251415 -> -8930111

# Row 64
$ws.Cells.Item(64, 8).Value = 2633.4048  # H64: 2621.0962 -> 2633.4048
$ws.Cells.Item(64, 9).Value = 2564.3572  # I64: 2617.261 -> 2564.3572
$ws.Cells.Item(64, 10).Value = 2771.5  # J64: 2624.138 -> 2771.5
$ws.Cells.Item(64, 11).Value = 2564.3572  # K64: 2617.261 -> 2564.3572
$ws.Cells.Item(64, 12).Value = 2771.5  # L64: 2624.138 -> 2771.5
$ws.Cells.Item(64, 13).Value = -2316.3572  # M64: -2369.261 -> -2316.3572
$ws.Cells.Item(64, 14).Value = -3267.5  # N64: -3120.138 -> -3267.5

# Row 67
$ws.Cells.Item(67, 8).Value = 2633.4048  # H67: 2621.0962 -> 2633.4048
$ws.Cells.Item(67, 9).Value = 2564.3572  # I67: 2617.261 -> 2564.3572
$ws.Cells.Item(67, 10).Value = 2771.5  # J67: 2624.138 -> 2771.5
$ws.Cells.Item(67, 11).Value = 2564.3572  # K67: 2617.261 -> 2564.3572
$ws.Cells.Item(67, 12).Value = 2771.5  # L67: 2624.138 -> 2771.5
$ws.Cells.Item(67, 13).Value = -1706.3572  # M67: -1759.261 -> -1706.3572
$ws.Cells.Item(67, 14).Value = -4487.5  # N67: -4340.138 -> -4487.5

# Row 106
$ws.Cells.Item(106, 8).Value = 2636.818  # H106: 2508.7083 -> 2636.818
$ws.Cells.Item(106, 9).Value = 2251  # I106: 2107.0625 -> 2251
$ws.Cells.Item(106, 11).Value = 2251  # K106: 2107.0625 -> 2251
$ws.Cells.Item(106, 13).Value = -1620  # M106: -1476.0625 -> -1620

# Row 132
$ws.Cells.Item(132, 8).Value = 7941412  # H132: 7697347.5 -> 7941412
$ws.Cells.Item(132, 9).Value = 9263401  # I132: 10208673 -> 9263401
$ws.Cells.Item(132, 10).Value = 9477.666999999999  # J132: 6414.75 -> 9477.666999999999
$ws.Cells.Item(132, 11).Value = 27790203  # K132: 30626019 -> 27790203
$ws.Cells.Item(132, 12).Value = 28433.001  # L132: 19244.25 -> 28433.001
$ws.Cells.Item(132, 13).Value = -27787673  # M132: -30623489 -> -27787673
$ws.Cells.Item(132, 14).Value = -33493.001  # N132: -24304.25 -> -33493.001

# Row 135
$ws.Cells.Item(135, 8).Value = 1065.2941  # H135: 764.24286 -> 1065.2941
$ws.Cells.Item(135, 9).Value = 904.1724  # I135: 599.371 -> 904.1724
$ws.Cells.Item(135, 10).Value = 1999.8  # J135: 2042 -> 1999.8
$ws.Cells.Item(135, 11).Value = 8137.551600000001  # K135: 5394.339 -> 8137.551600000001
$ws.Cells.Item(135, 12).Value = 17998.2  # L135: 18378 -> 17998.2
$ws.Cells.Item(135, 13).Value = -5602.551600000001  # M135: -2859.339 -> -5602.551600000001
$ws.Cells.Item(135, 14).Value = -23068.2  # N135: -23448 -> -23068.2

# Row 138
$ws.Cells.Item(138, 8).Value = 1573.5  # H138: 1666.2693 -> 1573.5
$ws.Cells.Item(138, 9).Value = 655.71875  # I138: 811.84375 -> 655.71875
$ws.Cells.Item(138, 10).Value = 3041.95  # J138: 3033.35 -> 3041.95
$ws.Cells.Item(138, 11).Value = 1967.15625  # K138: 2435.53125 -> 1967.15625
$ws.Cells.Item(138, 12).Value = 9125.849999999999  # L138: 9100.049999999999 -> 9125.849999999999
$ws.Cells.Item(138, 13).Value = 3172.84375  # M138: 2704.46875 -> 3172.84375
$ws.Cells.Item(138, 14).Value = -19405.85  # N138: -19380.05 -> -19405.85

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 627.5  # H2: 1070 -> 627.5
$ws.Cells.Item(2, 9).Value = 400.42856  # I2: 1023.3333 -> 400.42856
$ws.Cells.Item(2, 10).Value = 945.4  # J2: 1163.3334 -> 945.4
$ws.Cells.Item(2, 11).Value = 400.42856  # K2: 1023.3333 -> 400.42856
$ws.Cells.Item(2, 12).Value = 945.4  # L2: 1163.3334 -> 945.4
$ws.Cells.Item(2, 13).Value = -287.42856  # M2: -910.3333 -> -287.42856
$ws.Cells.Item(2, 14).Value = -1171.4  # N2: -1389.3334 -> -1171.4

# Row 8
$ws.Cells.Item(8, 8).Value = 652.5  # H8: 305 -> 652.5
$ws.Cells.Item(8, 9).Value = 652.5  # I8: 305 -> 652.5
$ws.Cells.Item(8, 11).Value = 652.5  # K8: 305 -> 652.5
$ws.Cells.Item(8, 13).Value = -508.5  # M8: -161 -> -508.5

# Row 32
$ws.Cells.Item(32, 8).Value = 1038.62  # H32: 797.14 -> 1038.62
$ws.Cells.Item(32, 9).Value = 998.96844  # I32: 760.14435 -> 998.96844
$ws.Cells.Item(32, 10).Value = 1792  # J32: 1993.3334 -> 1792
$ws.Cells.Item(32, 11).Value = 998.96844  # K32: 760.14435 -> 998.96844
$ws.Cells.Item(32, 12).Value = 1792  # L32: 1993.3334 -> 1792
$ws.Cells.Item(32, 13).Value = -711.96844  # M32: -473.14435 -> -711.96844
$ws.Cells.Item(32, 14).Value = -2366  # N32: -2567.3334 -> -2366

# Row 116
$ws.Cells.Item(116, 8).Value = 627.5  # H116: 1070 -> 627.5
$ws.Cells.Item(116, 9).Value = 400.42856  # I116: 1023.3333 -> 400.42856
$ws.Cells.Item(116, 10).Value = 945.4  # J116: 1163.3334 -> 945.4
$ws.Cells.Item(116, 11).Value = 400.42856  # K116: 1023.3333 -> 400.42856
$ws.Cells.Item(116, 12).Value = 945.4  # L116: 1163.3334 -> 945.4
$ws.Cells.Item(116, 13).Value = 1893.57144  # M116: 1270.6667 -> 1893.57144
$ws.Cells.Item(116, 14).Value = -5533.4  # N116: -5751.3334 -> -5533.4

# Row 122
$ws.Cells.Item(122, 8).Value = 1383  # H122: 840 -> 1383
$ws.Cells.Item(122, 9).Value = 1199.5  # I122: 840 -> 1199.5
$ws.Cells.Item(122, 10).Value = 1750  # J122: 0 -> 1750
$ws.Cells.Item(122, 11).Value = 3598.5  # K122: 2520 -> 3598.5
$ws.Cells.Item(122, 12).Value = 5250  # L122: 0 -> 5250
$ws.Cells.Item(122, 13).Value = -1148.5  # M122: -70 -> -1148.5
$ws.Cells.Item(122, 14).Value = -10150  # N122: None -> -10150

# Row 132
$ws.Cells.Item(132, 8).Value = 1899456.1  # H132: 1402050.5 -> 1899456.1
$ws.Cells.Item(132, 9).Value = 1427.6  # I132: 1110.5714 -> 1427.6
$ws.Cells.Item(132, 10).Value = 3678858  # J132: 4203930.5 -> 3678858
$ws.Cells.Item(132, 11).Value = 4282.799999999999  # K132: 3331.7142 -> 4282.799999999999
$ws.Cells.Item(132, 12).Value = 11036574  # L132: 12611791.5 -> 11036574
$ws.Cells.Item(132, 13).Value = -1752.799999999999  # M132: -801.7142000000003 -> -1752.799999999999
$ws.Cells.Item(132, 14).Value = -11041634  # N132: -12616851.5 -> -11041634

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 627.5  # H3: 1070 -> 627.5
$ws.Cells.Item(3, 9).Value = 400.42856  # I3: 1023.3333 -> 400.42856
$ws.Cells.Item(3, 10).Value = 945.4  # J3: 1163.3334 -> 945.4
$ws.Cells.Item(3, 11).Value = 400.42856  # K3: 1023.3333 -> 400.42856
$ws.Cells.Item(3, 12).Value = 945.4  # L3: 1163.3334 -> 945.4
$ws.Cells.Item(3, 13).Value = -286.42856  # M3: -909.3333 -> -286.42856
$ws.Cells.Item(3, 14).Value = -1173.4  # N3: -1391.3334 -> -1173.4

# Row 20
$ws.Cells.Item(20, 8).Value = 3490.6365  # H20: 4074.3333 -> 3490.6365
$ws.Cells.Item(20, 9).Value = 3198.5  # I20: 4165 -> 3198.5
$ws.Cells.Item(20, 10).Value = 4269.6665  # J20: 4029 -> 4269.6665
$ws.Cells.Item(20, 11).Value = 3198.5  # K20: 4165 -> 3198.5
$ws.Cells.Item(20, 12).Value = 4269.6665  # L20: 4029 -> 4269.6665
$ws.Cells.Item(20, 13).Value = -2951.5  # M20: -3918 -> -2951.5
$ws.Cells.Item(20, 14).Value = -4763.6665  # N20: -4523 -> -4763.6665

# Row 30
$ws.Cells.Item(30, 8).Value = 25005.5  # H30: 32666.666 -> 25005.5
$ws.Cells.Item(30, 9).Value = 0  # I30: 3000 -> 0
$ws.Cells.Item(30, 10).Value = 25005.5  # J30: 47500 -> 25005.5
$ws.Cells.Item(30, 11).Value = 0  # K30: 3000 -> 0
$ws.Cells.Item(30, 12).Value = 25005.5  # L30: 47500 -> 25005.5
$ws.Cells.Item(30, 13).ClearContents()  # M30: -2875 -> (blank)
$ws.Cells.Item(30, 14).Value = -25255.5  # N30: -47750 -> -25255.5

# Row 86
$ws.Cells.Item(86, 8).Value = 1108963.6  # H86: 1164558 -> 1108963.6
$ws.Cells.Item(86, 9).Value = 1393  # I86: 2480.125 -> 1393
$ws.Cells.Item(86, 10).Value = 2327291.2  # J86: 1939276.6 -> 2327291.2
$ws.Cells.Item(86, 11).Value = 1393  # K86: 2480.125 -> 1393
$ws.Cells.Item(86, 12).Value = 2327291.2  # L86: 1939276.6 -> 2327291.2
$ws.Cells.Item(86, 13).Value = -270  # M86: -1357.125 -> -270
$ws.Cells.Item(86, 14).Value = -2329537.2  # N86: -1941522.6 -> -2329537.2

# Row 89
$ws.Cells.Item(89, 8).Value = 1108963.6  # H89: 1164558 -> 1108963.6
$ws.Cells.Item(89, 9).Value = 1393  # I89: 2480.125 -> 1393
$ws.Cells.Item(89, 10).Value = 2327291.2  # J89: 1939276.6 -> 2327291.2
$ws.Cells.Item(89, 11).Value = 6965  # K89: 12400.625 -> 6965
$ws.Cells.Item(89, 12).Value = 11636456  # L89: 9696383 -> 11636456
$ws.Cells.Item(89, 13).Value = -1349  # M89: -6784.625 -> -1349
$ws.Cells.Item(89, 14).Value = -11647688  # N89: -9707615 -> -11647688

# Row 107
$ws.Cells.Item(107, 8).Value = 35715036  # H107: 45455390 -> 35715036
$ws.Cells.Item(107, 9).Value = 38462116  # I107: 50000630 -> 38462116
$ws.Cells.Item(107, 11).Value = 38462116  # K107: 50000630 -> 38462116
$ws.Cells.Item(107, 13).Value = -38460196  # M107: -49998710 -> -38460196

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 28943.143  # H22: 10282.7 -> 28943.143
$ws.Cells.Item(22, 9).Value = 25153  # I22: 233.28572 -> 25153
$ws.Cells.Item(22, 10).Value = 33996.668  # J22: 33731.332 -> 33996.668
$ws.Cells.Item(22, 11).Value = 25153  # K22: 233.28572 -> 25153
$ws.Cells.Item(22, 12).Value = 33996.668  # L22: 33731.332 -> 33996.668
$ws.Cells.Item(22, 13).Value = -24803  # M22: 116.71428 -> -24803
$ws.Cells.Item(22, 14).Value = -34696.668  # N22: -34431.332 -> -34696.668

# Row 25
$ws.Cells.Item(25, 8).Value = 1300  # H25: 1705.2 -> 1300
$ws.Cells.Item(25, 9).Value = 1300  # I25: 631.5 -> 1300
$ws.Cells.Item(25, 10).Value = 0  # J25: 6000 -> 0
$ws.Cells.Item(25, 11).Value = 1300  # K25: 631.5 -> 1300
$ws.Cells.Item(25, 12).Value = 0  # L25: 6000 -> 0
$ws.Cells.Item(25, 13).Value = -1126  # M25: -457.5 -> -1126
$ws.Cells.Item(25, 14).ClearContents()  # N25: -6348 -> (blank)

# Row 41
$ws.Cells.Item(41, 8).Value = 0  # H41: 5535 -> 0
$ws.Cells.Item(41, 9).Value = 0  # I41: 5535 -> 0
$ws.Cells.Item(41, 11).Value = 0  # K41: 5535 -> 0
$ws.Cells.Item(41, 13).ClearContents()  # M41: -5107 -> (blank)

# Row 134
$ws.Cells.Item(134, 8).Value = 903.6579  # H134: 1161.3214 -> 903.6579
$ws.Cells.Item(134, 9).Value = 894.35486  # I134: 1228.7142 -> 894.35486
$ws.Cells.Item(134, 10).Value = 944.8570999999999  # J134: 959.1429000000001 -> 944.8570999999999
$ws.Cells.Item(134, 11).Value = 2683.06458  # K134: 3686.1426 -> 2683.06458
$ws.Cells.Item(134, 12).Value = 2834.5713  # L134: 2877.4287 -> 2834.5713
$ws.Cells.Item(134, 13).Value = -148.0645800000002  # M134: -1151.1426 -> -148.0645800000002
$ws.Cells.Item(134, 14).Value = -7904.5713  # N134: -7947.4287 -> -7904.5713

$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Cells.Item(92, 8).Value = 6920.294  # H92: 9708.333000000001 -> 6920.294
$ws.Cells.Item(92, 9).Value = 661.5  # I92: 266.66666 -> 661.5
$ws.Cells.Item(92, 10).Value = 8846.076999999999  # J92: 12855.556 -> 8846.076999999999
$ws.Cells.Item(92, 11).Value = 1984.5  # K92: 799.9999799999999 -> 1984.5
$ws.Cells.Item(92, 12).Value = 26538.231  # L92: 38566.66800000001 -> 26538.231
$ws.Cells.Item(92, 13).Value = -736.5  # M92: 448.0000200000001 -> -736.5
$ws.Cells.Item(92, 14).Value = -29034.231  # N92: -41062.66800000001 -> -29034.231

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 380.8  # H3: 1000 -> 380.8
$ws.Cells.Item(3, 9).Value = 351  # I3: 1000 -> 351
$ws.Cells.Item(3, 10).Value = 500  # J3: 0 -> 500
$ws.Cells.Item(3, 11).Value = 351  # K3: 1000 -> 351
$ws.Cells.Item(3, 12).Value = 500  # L3: 0 -> 500
$ws.Cells.Item(3, 13).Value = -235  # M3: -884 -> -235
$ws.Cells.Item(3, 14).Value = -732  # N3: None -> -732

# Row 70
$ws.Cells.Item(70, 8).Value = 5348.778  # H70: 5197.0312 -> 5348.778
$ws.Cells.Item(70, 9).Value = 4588.2354  # I70: 4483.6313 -> 4588.2354
$ws.Cells.Item(70, 10).Value = 6641.7  # J70: 6239.6924 -> 6641.7
$ws.Cells.Item(70, 11).Value = 4588.2354  # K70: 4483.6313 -> 4588.2354
$ws.Cells.Item(70, 12).Value = 6641.7  # L70: 6239.6924 -> 6641.7
$ws.Cells.Item(70, 13).Value = -4318.2354  # M70: -4213.6313 -> -4318.2354
$ws.Cells.Item(70, 14).Value = -7181.7  # N70: -6779.6924 -> -7181.7

# Row 73
$ws.Cells.Item(73, 8).Value = 5348.778  # H73: 5197.0312 -> 5348.778
$ws.Cells.Item(73, 9).Value = 4588.2354  # I73: 4483.6313 -> 4588.2354
$ws.Cells.Item(73, 10).Value = 6641.7  # J73: 6239.6924 -> 6641.7
$ws.Cells.Item(73, 11).Value = 4588.2354  # K73: 4483.6313 -> 4588.2354
$ws.Cells.Item(73, 12).Value = 6641.7  # L73: 6239.6924 -> 6641.7
$ws.Cells.Item(73, 13).Value = -3652.2354  # M73: -3547.6313 -> -3652.2354
$ws.Cells.Item(73, 14).Value = -8513.700000000001  # N73: -8111.6924 -> -8513.700000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 1475.5  # H68: 1500 -> 1475.5
$ws.Cells.Item(68, 9).Value = 1451  # I68: 0 -> 1451
$ws.Cells.Item(68, 11).Value = 1451  # K68: 0 -> 1451
$ws.Cells.Item(68, 13).Value = -702  # M68: None -> -702

# Row 71
$ws.Cells.Item(71, 8).Value = 1475.5  # H71: 1500 -> 1475.5
$ws.Cells.Item(71, 9).Value = 1451  # I71: 0 -> 1451
$ws.Cells.Item(71, 11).Value = 7255  # K71: 0 -> 7255
$ws.Cells.Item(71, 13).Value = -3511  # M71: None -> -3511

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 1157.7307  # H126: 1113 -> 1157.7307
$ws.Cells.Item(126, 9).Value = 931.86365  # I126: 888.5599999999999 -> 931.86365
$ws.Cells.Item(126, 10).Value = 2400  # J126: 2983.3333 -> 2400
$ws.Cells.Item(126, 11).Value = 2795.59095  # K126: 2665.68 -> 2795.59095
$ws.Cells.Item(126, 12).Value = 7200  # L126: 8949.999899999999 -> 7200
$ws.Cells.Item(126, 13).Value = -325.5909499999998  # M126: -195.6799999999998 -> -325.5909499999998
$ws.Cells.Item(126, 14).Value = -12140  # N126: -13889.9999 -> -12140

